$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update row 10 (Objetivos:) B/C text: replace misplaced professor name
#        with the actual Portuguese objectives paragraph.
$ws.Range("B10").Value = "Abordar os princípios e conceitos da evolução biológica e da ecologia em suas diferentes escalas: populações, comunidades e ecossistemas, dentro do enfoque da sustentabilidade. Apresentar aspectos econômicos, sociais e culturais da sociedade envolvidos na preservação da biodiversidade e dos ecossistemas."
$ws.Range("C10").Value = "Abordar os princípios e conceitos da evolução biológica e da ecologia em suas diferentes escalas: populações, comunidades e ecossistemas, dentro do enfoque da sustentabilidade. Apresentar aspectos econômicos, sociais e culturais da sociedade envolvidos na preservação da biodiversidade e dos ecossistemas."

# --- 2) Insert a new row after row 12 ("Docentes responsáveis:") to hold the
#        professor name in its own row (B13/C13), pushing everything else down.
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()
# Copy the B/C column formatting (wrap text, top-valign, red font for C) from
# an existing data row instead of rebuilding it by hand, so no stray/unused
# cell-format entries get added to the style table.
$ws.Range("B11").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C13").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"

# --- 3) Row 14 ("Programa resumido:") B/C: was "Semestral", now the Portuguese
#        short-syllabus paragraph.
$ws.Range("B14").Value = "1. Problemas ambientais, causas e soluções2. Ecossistemas: o que são e como funcionam3. Biodiversidade e evolução4. Biodiversidade, interações de espécies e controle da população5. A população humana e seu impacto"
$ws.Range("C14").Value = "1. Problemas ambientais, causas e soluções2. Ecossistemas: o que são e como funcionam3. Biodiversidade e evolução4. Biodiversidade, interações de espécies e controle da população5. A população humana e seu impacto"

# --- 4) Row 16 ("Programa:") B/C: was the (wrong) activation date, now the
#        same Portuguese syllabus paragraph as row 14.
$ws.Range("B16").Value = "1. Problemas ambientais, causas e soluções2. Ecossistemas: o que são e como funcionam3. Biodiversidade e evolução4. Biodiversidade, interações de espécies e controle da população5. A população humana e seu impacto"
$ws.Range("C16").Value = "1. Problemas ambientais, causas e soluções2. Ecossistemas: o que são e como funcionam3. Biodiversidade e evolução4. Biodiversidade, interações de espécies e controle da população5. A população humana e seu impacto"

# --- 5) Rows 19/20/21 ("Método:"/"Critério:"/"Norma de recuperação:") B/C:
#        each now holds the text that used to belong to the row above it
#        (the whole evaluation block shifted down by one row).
$ws.Range("B19").Value = "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e seminários."
$ws.Range("C19").Value = "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e seminários."

$ws.Range("B20").Value = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Range("C20").Value = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."

$ws.Range("B21").Value = "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2."
$ws.Range("C21").Value = "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2."

# --- 6) New row 22 at the end: Bibliografia. Again, copy B/C formatting from
#        an existing "big text" row (A/B/C styles only, not row height).
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Rows.Item(22).RowHeight = 120
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Básica:MILLER, G.T.; SPOOLMAN, S.E. 2012. Ecologia e sustentabilidade. Cengage Learning. 412p.Complementar:BEGON, M., J.L. HARPER & C.R. TOWNSEND. 2005. Ecology. From Individuals to Communities. Blackwell Science.RICKLEFS, R.E. 2003. A economia da natureza. Guanabara Koogan.RICKLEFS, R.E. & G.L. MILLER. 2000. Ecology. W.H. Freeman and Co.TOWNSEND, C.R., M. BEGON. & J.L. HARPER 2006. Fundamentos em ecologia. Artmed."
$ws.Range("C22").Value = "Básica:MILLER, G.T.; SPOOLMAN, S.E. 2012. Ecologia e sustentabilidade. Cengage Learning. 412p.Complementar:BEGON, M., J.L. HARPER & C.R. TOWNSEND. 2005. Ecology. From Individuals to Communities. Blackwell Science.RICKLEFS, R.E. 2003. A economia da natureza. Guanabara Koogan.RICKLEFS, R.E. & G.L. MILLER. 2000. Ecology. W.H. Freeman and Co.TOWNSEND, C.R., M. BEGON. & J.L. HARPER 2006. Fundamentos em ecologia. Artmed."
